$d = $word.ActiveDocument

# Remove the stray "}}" left over between "...cannot know in advance" and the
# following "{% if situation == ..." template tag.
$d.Content.Find.Execute(
    "cannot know in advance}}{% if situation",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "cannot know in advance{% if situation",
    2
)
